# Generate Report for Handoff
# A new source file (47595699-2b4b-440b-a5e4-6bd920b4c599) replaces the
# previously-tracked one, and a second new file
# (ffffcd80ed84-1304-4fc4-a090-4ff269d61d9f) is now also being handed off.
# Every worksheet gains one extra data row for it.

$wb = $excel.ActiveWorkbook

$oldMd   = "7bfd46b1-2996-4efb-aa21-5be5a8f9f7f5.md"
$newMd1  = "47595699-2b4b-440b-a5e4-6bd920b4c599.md"
$newMd2  = "ffffcd80ed84-1304-4fc4-a090-4ff269d61d9f.md"

$oldHash = "0c9fce0f97ddda2e389cf7169b73745d5315e6e4"
$newHash = "8efaaefef328ff2296b58269c3d3b5a70aadfaf6"

$zhXlf = "47595699-2b4b-440b-a5e4-6bd920b4c599." + $newHash + ".zh-cn.xlf"
$deXlf = "47595699-2b4b-440b-a5e4-6bd920b4c599." + $newHash + ".de-de.xlf"

$zhHandoffDatetime = "2016-02-23 09:16:45"
$deHandoffDatetime = "2016-02-23 09:16:57"

$mdBase = "https://github.com/OpenLocalizationTest/oltest/blob/cff02456aa3344f8dc914e42029d0b941c73b8dc/e2e/"
$cfgUrl = "https://github.com/OpenLocalizationTest/oltest/blob/cff02456aa3344f8dc914e42029d0b941c73b8dc/.localization-config"
$zhXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/d9529e0f3df2b2554a5be06d8d6eb91c17153be5/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/xinjiang/ht/" + $zhXlf
$deXlfUrl = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/6676a06baa62f97860a615ba466d3afac2671208/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/xinjiang/ht/" + $deXlf

# ---------------------------------------------------------------------
# Overview sheet: File Name / zh-cn / de-de summary
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Hyperlinks.Delete()
$wsOverview.Rows("3:3").Insert()

$wsOverview.Range("A2").Value = $newMd1
$wsOverview.Range("A3").Value = $newMd2
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"
$wsOverview.Range("A4").Value = ".localization-config"
$wsOverview.Range("B4").Value = "Not to be localized"
$wsOverview.Range("C4").Value = "Not to be localized"

$wsOverview.Hyperlinks.Add($wsOverview.Range("A2"), $mdBase + $newMd1, "", "", $newMd1)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A3"), $mdBase + $newMd2, "", "", $newMd2)
$wsOverview.Hyperlinks.Add($wsOverview.Range("A4"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# zh-cn sheet
# ---------------------------------------------------------------------
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Hyperlinks.Delete()
$wsZh.Rows("3:3").Insert()

$wsZh.Range("A2").Value = $newMd1
$wsZh.Range("C2").Value = $zhXlf
$wsZh.Range("D2").Value = $zhHandoffDatetime

$wsZh.Range("A3").Value = $newMd2
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = $zhXlf
$wsZh.Range("C3").NumberFormat = "General"
$wsZh.Range("D3").Value = $zhHandoffDatetime
$wsZh.Range("G3").Value = "0001-01-01 00:00:00"
$wsZh.Range("H3").Value = "Include"

$wsZh.Range("A4").Value = ".localization-config"
$wsZh.Range("B4").Value = "Not to be localized"
$wsZh.Range("D4").Value = "0001-01-01 00:00:00"
$wsZh.Range("G4").Value = "0001-01-01 00:00:00"
$wsZh.Range("H4").Value = "Ignored"

$wsZh.Hyperlinks.Add($wsZh.Range("A2"), $mdBase + $newMd1, "", "", $newMd1)
$wsZh.Hyperlinks.Add($wsZh.Range("C2"), $zhXlfUrl, "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A3"), $mdBase + $newMd2, "", "", $newMd2)
$wsZh.Hyperlinks.Add($wsZh.Range("C3"), $zhXlfUrl, "", "", $zhXlf)
$wsZh.Hyperlinks.Add($wsZh.Range("A4"), $cfgUrl, "", "", ".localization-config")

# ---------------------------------------------------------------------
# de-de sheet
# ---------------------------------------------------------------------
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Hyperlinks.Delete()
$wsDe.Rows("3:3").Insert()

$wsDe.Range("A2").Value = $newMd1
$wsDe.Range("C2").Value = $deXlf
$wsDe.Range("D2").Value = $deHandoffDatetime

$wsDe.Range("A3").Value = $newMd2
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = $deXlf
$wsDe.Range("C3").NumberFormat = "General"
$wsDe.Range("D3").Value = $deHandoffDatetime
$wsDe.Range("G3").Value = "0001-01-01 00:00:00"
$wsDe.Range("H3").Value = "Include"

$wsDe.Range("A4").Value = ".localization-config"
$wsDe.Range("B4").Value = "Not to be localized"
$wsDe.Range("D4").Value = "0001-01-01 00:00:00"
$wsDe.Range("G4").Value = "0001-01-01 00:00:00"
$wsDe.Range("H4").Value = "Ignored"

$wsDe.Hyperlinks.Add($wsDe.Range("A2"), $mdBase + $newMd1, "", "", $newMd1)
$wsDe.Hyperlinks.Add($wsDe.Range("C2"), $deXlfUrl, "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A3"), $mdBase + $newMd2, "", "", $newMd2)
$wsDe.Hyperlinks.Add($wsDe.Range("C3"), $deXlfUrl, "", "", $deXlf)
$wsDe.Hyperlinks.Add($wsDe.Range("A4"), $cfgUrl, "", "", ".localization-config")
